$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("J2").Value = 2229
$ws.Range("J3").Value = 2335
$ws.Range("J4").Value = 526
$ws.Range("J5").Value = 169
$ws.Range("J6").Value = 2937
$ws.Range("J7").Value = 8196

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("J6").Value = 28
$ws.Range("J7").Value = 85

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("J2").Value = 80
$ws.Range("J3").Value = 101
$ws.Range("J7").Value = 276

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("J2").Value = 24
$ws.Range("J3").Value = 16
$ws.Range("J6").Value = 13
$ws.Range("J7").Value = 62

$ws = $wb.Worksheets.Item("New City")
$ws.Range("J3").Value = 60
$ws.Range("J4").Value = 12
$ws.Range("J6").Value = 78
$ws.Range("J7").Value = 218

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J4").Value = 34
$ws.Range("J7").Value = 249
$ws.Range("J8").Value = 518
$ws.Range("J10").Value = 49
$ws.Range("J20").Value = 170
$ws.Range("J25").Value = 48
$ws.Range("J27").Value = 49
$ws.Range("J28").Value = 3
$ws.Range("J29").Value = 461
$ws.Range("J31").Value = 62
$ws.Range("J32").Value = 15
$ws.Range("J33").Value = 338
$ws.Range("J37").Value = 276
$ws.Range("J42").Value = 311
$ws.Range("J43").Value = 79
$ws.Range("J46").Value = 27
$ws.Range("J47").Value = 73
$ws.Range("J48").Value = 76
$ws.Range("J52").Value = 199
$ws.Range("J53").Value = 79
$ws.Range("J54").Value = 167
$ws.Range("J63").Value = 29
$ws.Range("J64").Value = 55
$ws.Range("J65").Value = 218
$ws.Range("J66").Value = 20
$ws.Range("J77").Value = 61
$ws.Range("J78").Value = 113
$ws.Range("J79").Value = 248
$ws.Range("J83").Value = 194
$ws.Range("J85").Value = 389
$ws.Range("J89").Value = 85
$ws.Range("J90").Value = 91
$ws.Range("J91").Value = 92
$ws.Range("J94").Value = 65
$ws.Range("J95").Value = 123
$ws.Range("J101").Value = 8196

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("J6").Value = 58
$ws.Range("J7").Value = 194

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("J2").Value = 43
$ws.Range("J3").Value = 37
$ws.Range("J7").Value = 123

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("J2").Value = 88
$ws.Range("J7").Value = 338

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("J2").Value = 41
$ws.Range("J7").Value = 167

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("J2").Value = 135
$ws.Range("J3").Value = 159
$ws.Range("J4").Value = 25
$ws.Range("J6").Value = 126
$ws.Range("J7").Value = 461

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("J6").Value = 38
$ws.Range("J7").Value = 76

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("J2").Value = 94
$ws.Range("J6").Value = 113
$ws.Range("J7").Value = 389

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("J3").Value = 66
$ws.Range("J7").Value = 311

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("J6").Value = 20
$ws.Range("J7").Value = 49

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("J3").Value = 39
$ws.Range("J5").Value = 3
$ws.Range("J7").Value = 113

$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("J2").Value = 8
$ws.Range("J7").Value = 27

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("J6").Value = 17
$ws.Range("J7").Value = 92

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("J3").Value = 90
$ws.Range("J6").Value = 69
$ws.Range("J7").Value = 248

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("J2").Value = 16
$ws.Range("J6").Value = 19
$ws.Range("J7").Value = 55

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("J6").Value = 45
$ws.Range("J7").Value = 170

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("J3").Value = 59
$ws.Range("J7").Value = 199

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("J6").Value = 39
$ws.Range("J7").Value = 65

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("J6").Value = 11
$ws.Range("J7").Value = 48

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("J2").Value = 20
$ws.Range("J7").Value = 73

$ws = $wb.Worksheets.Item("North Center")
$ws.Range("J6").Value = 12
$ws.Range("J7").Value = 20

$ws = $wb.Worksheets.Item("Galewood")
$ws.Range("J3").Value = 4
$ws.Range("J7").Value = 15

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("J2").Value = 162
$ws.Range("J3").Value = 169
$ws.Range("J4").Value = 24
$ws.Range("J6").Value = 149
$ws.Range("J7").Value = 518

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("J2").Value = 12
$ws.Range("J7").Value = 49

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("J2").Value = 27
$ws.Range("J3").Value = 25
$ws.Range("J7").Value = 91

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("J3").Value = 15
$ws.Range("J7").Value = 79

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("J2").Value = 14
$ws.Range("J6").Value = 45
$ws.Range("J7").Value = 79

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("J2").Value = 18
$ws.Range("J7").Value = 61

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("J2").Value = 80
$ws.Range("J3").Value = 75
$ws.Range("J7").Value = 249

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("J3").Value = 8
$ws.Range("J7").Value = 34

$ws = $wb.Worksheets.Item("Edison Park")
$ws.Range("J6").Value = 1
$ws.Range("J7").Value = 3
